$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 12).Value = "stimuli/img_5il0t.png"
$ws.Cells.Item(2, 13).Value = 48.09523809523809
$ws.Cells.Item(2, 14).Value = 30.90476190476191
$ws.Cells.Item(2, 15).Value = 39.5
$ws.Cells.Item(2, 16).Value = 42
$ws.Cells.Item(2, 17).Value = 2
$ws.Cells.Item(2, 18).Value = 2
$ws.Cells.Item(2, 19).Value = 2
$ws.Cells.Item(2, 20).Value = 2
$ws.Cells.Item(2, 21).Value = 2
$ws.Cells.Item(2, 22).Value = 2

# Row 3
$ws.Cells.Item(3, 12).Value = "stimuli/img_scrdm.png"
$ws.Cells.Item(3, 13).Value = 78.675
$ws.Cells.Item(3, 14).Value = 57.9
$ws.Cells.Item(3, 15).Value = 68.28749999999999
$ws.Cells.Item(3, 17).Value = 7
$ws.Cells.Item(3, 18).Value = 7
$ws.Cells.Item(3, 19).Value = 7
$ws.Cells.Item(3, 20).Value = 7
$ws.Cells.Item(3, 21).Value = 7
$ws.Cells.Item(3, 22).Value = 7

# Row 4
$ws.Cells.Item(4, 12).Value = "stimuli/img_wijef.png"
$ws.Cells.Item(4, 13).Value = 69.875
$ws.Cells.Item(4, 14).Value = 48.025
$ws.Cells.Item(4, 15).Value = 58.95
$ws.Cells.Item(4, 16).Value = 40
$ws.Cells.Item(4, 17).Value = 5
$ws.Cells.Item(4, 18).Value = 5
$ws.Cells.Item(4, 19).Value = 5
$ws.Cells.Item(4, 20).Value = 5
$ws.Cells.Item(4, 21).Value = 5
$ws.Cells.Item(4, 22).Value = 5

# Row 5
$ws.Cells.Item(5, 9).ClearContents()
$ws.Cells.Item(5, 10).Value = "new"
$ws.Cells.Item(5, 11).Value = "f"
$ws.Cells.Item(5, 12).Value = "stimuli/img_4wq98.png"
$ws.Cells.Item(5, 13).Value = 78.48387096774194
$ws.Cells.Item(5, 14).Value = 58.12903225806452
$ws.Cells.Item(5, 15).Value = 68.30645161290323
$ws.Cells.Item(5, 16).Value = 31
$ws.Cells.Item(5, 17).Value = 7
$ws.Cells.Item(5, 18).Value = 7
$ws.Cells.Item(5, 19).Value = 7
$ws.Cells.Item(5, 20).Value = 7
$ws.Cells.Item(5, 21).Value = 7
$ws.Cells.Item(5, 22).Value = 7

# Row 6
$ws.Cells.Item(6, 9).Value = "target"
$ws.Cells.Item(6, 10).Value = "old"
$ws.Cells.Item(6, 11).Value = "j"
$ws.Cells.Item(6, 12).Value = "stimuli/img_72fmj.png"
$ws.Cells.Item(6, 13).Value = 53.87179487179487
$ws.Cells.Item(6, 14).Value = 36.02564102564103
$ws.Cells.Item(6, 15).Value = 44.94871794871795
$ws.Cells.Item(6, 16).Value = 39
$ws.Cells.Item(6, 17).Value = 3
$ws.Cells.Item(6, 18).Value = 3
$ws.Cells.Item(6, 19).Value = 3
$ws.Cells.Item(6, 20).Value = 3
$ws.Cells.Item(6, 21).Value = 3
$ws.Cells.Item(6, 22).Value = 3

# Row 7
$ws.Cells.Item(7, 12).Value = "stimuli/img_zi682.png"
$ws.Cells.Item(7, 13).Value = 84.59999999999999
$ws.Cells.Item(7, 14).Value = 69.52500000000001
$ws.Cells.Item(7, 15).Value = 77.0625
$ws.Cells.Item(7, 16).Value = 40
$ws.Cells.Item(7, 17).Value = 9
$ws.Cells.Item(7, 18).Value = 9
$ws.Cells.Item(7, 19).Value = 9
$ws.Cells.Item(7, 20).Value = 9
$ws.Cells.Item(7, 21).Value = 9
$ws.Cells.Item(7, 22).Value = 9

# Row 9
$ws.Cells.Item(9, 9).ClearContents()
$ws.Cells.Item(9, 10).Value = "new"
$ws.Cells.Item(9, 11).Value = "f"
$ws.Cells.Item(9, 12).Value = "stimuli/img_ce55l.png"
$ws.Cells.Item(9, 13).Value = 82.23809523809524
$ws.Cells.Item(9, 14).Value = 64.07142857142857
$ws.Cells.Item(9, 15).Value = 73.1547619047619
$ws.Cells.Item(9, 17).Value = 8
$ws.Cells.Item(9, 18).Value = 8
$ws.Cells.Item(9, 19).Value = 8
$ws.Cells.Item(9, 20).Value = 8
$ws.Cells.Item(9, 21).Value = 8
$ws.Cells.Item(9, 22).Value = 8

# Row 10
$ws.Cells.Item(10, 9).ClearContents()
$ws.Cells.Item(10, 10).Value = "new"
$ws.Cells.Item(10, 11).Value = "f"
$ws.Cells.Item(10, 12).Value = "stimuli/img_5yhyk.png"
$ws.Cells.Item(10, 13).Value = 46.375
$ws.Cells.Item(10, 14).Value = 31.325
$ws.Cells.Item(10, 15).Value = 38.85
$ws.Cells.Item(10, 16).Value = 40
$ws.Cells.Item(10, 17).Value = 2
$ws.Cells.Item(10, 18).Value = 2
$ws.Cells.Item(10, 19).Value = 2
$ws.Cells.Item(10, 20).Value = 2
$ws.Cells.Item(10, 21).Value = 2
$ws.Cells.Item(10, 22).Value = 2

# Row 11
$ws.Cells.Item(11, 9).ClearContents()
$ws.Cells.Item(11, 10).Value = "new"
$ws.Cells.Item(11, 11).Value = "f"
$ws.Cells.Item(11, 12).Value = "stimuli/img_zgg62.png"
$ws.Cells.Item(11, 13).Value = 82.18421052631579
$ws.Cells.Item(11, 14).Value = 63.52631578947368
$ws.Cells.Item(11, 15).Value = 72.85526315789474
$ws.Cells.Item(11, 16).Value = 38
$ws.Cells.Item(11, 17).Value = 8
$ws.Cells.Item(11, 18).Value = 8
$ws.Cells.Item(11, 19).Value = 8
$ws.Cells.Item(11, 20).Value = 8
$ws.Cells.Item(11, 21).Value = 8
$ws.Cells.Item(11, 22).Value = 8

# Row 12
$ws.Cells.Item(12, 12).Value = "stimuli/img_psgf7.png"
$ws.Cells.Item(12, 13).Value = 26
$ws.Cells.Item(12, 14).Value = 11.66666666666667
$ws.Cells.Item(12, 15).Value = 18.83333333333333
$ws.Cells.Item(12, 16).Value = 36
$ws.Cells.Item(12, 17).Value = 1
$ws.Cells.Item(12, 18).Value = 1
$ws.Cells.Item(12, 19).Value = 1
$ws.Cells.Item(12, 20).Value = 1
$ws.Cells.Item(12, 21).Value = 1
$ws.Cells.Item(12, 22).Value = 1

# Row 13
$ws.Cells.Item(13, 9).Value = "target"
$ws.Cells.Item(13, 10).Value = "old"
$ws.Cells.Item(13, 11).Value = "j"
$ws.Cells.Item(13, 12).Value = "stimuli/img_rvssl.png"
$ws.Cells.Item(13, 13).Value = 74.25
$ws.Cells.Item(13, 14).Value = 54.33333333333334
$ws.Cells.Item(13, 15).Value = 64.29166666666667
$ws.Cells.Item(13, 16).Value = 36
$ws.Cells.Item(13, 17).Value = 6
$ws.Cells.Item(13, 18).Value = 6
$ws.Cells.Item(13, 19).Value = 6
$ws.Cells.Item(13, 20).Value = 6
$ws.Cells.Item(13, 21).Value = 6
$ws.Cells.Item(13, 22).Value = 6

# Row 14
$ws.Cells.Item(14, 9).Value = "target"
$ws.Cells.Item(14, 10).Value = "old"
$ws.Cells.Item(14, 11).Value = "j"
$ws.Cells.Item(14, 12).Value = "stimuli/img_le8uf.png"
$ws.Cells.Item(14, 13).Value = 12.88888888888889
$ws.Cells.Item(14, 14).Value = 9.222222222222221
$ws.Cells.Item(14, 15).Value = 11.05555555555556
$ws.Cells.Item(14, 16).Value = 36
$ws.Cells.Item(14, 17).Value = 1
$ws.Cells.Item(14, 18).Value = 1
$ws.Cells.Item(14, 19).Value = 1
$ws.Cells.Item(14, 20).Value = 1
$ws.Cells.Item(14, 21).Value = 1
$ws.Cells.Item(14, 22).Value = 1

# Row 15
$ws.Cells.Item(15, 9).Value = "target"
$ws.Cells.Item(15, 10).Value = "old"
$ws.Cells.Item(15, 11).Value = "j"
$ws.Cells.Item(15, 12).Value = "stimuli/img_cmyvx.png"
$ws.Cells.Item(15, 13).Value = 64.25
$ws.Cells.Item(15, 14).Value = 40.09375
$ws.Cells.Item(15, 15).Value = 52.171875
$ws.Cells.Item(15, 16).Value = 32
$ws.Cells.Item(15, 17).Value = 4
$ws.Cells.Item(15, 18).Value = 4
$ws.Cells.Item(15, 19).Value = 4
$ws.Cells.Item(15, 20).Value = 4
$ws.Cells.Item(15, 21).Value = 4
$ws.Cells.Item(15, 22).Value = 4

# Row 16
$ws.Cells.Item(16, 9).Value = "target"
$ws.Cells.Item(16, 10).Value = "old"
$ws.Cells.Item(16, 11).Value = "j"
$ws.Cells.Item(16, 12).Value = "stimuli/img_f4jxo.png"
$ws.Cells.Item(16, 13).Value = 82.91666666666667
$ws.Cells.Item(16, 14).Value = 65.52777777777777
$ws.Cells.Item(16, 15).Value = 74.22222222222223
$ws.Cells.Item(16, 16).Value = 36
$ws.Cells.Item(16, 17).Value = 8
$ws.Cells.Item(16, 18).Value = 8
$ws.Cells.Item(16, 19).Value = 8
$ws.Cells.Item(16, 20).Value = 8
$ws.Cells.Item(16, 21).Value = 8
$ws.Cells.Item(16, 22).Value = 8

# Row 17
$ws.Cells.Item(17, 9).Value = "target"
$ws.Cells.Item(17, 10).Value = "old"
$ws.Cells.Item(17, 11).Value = "j"
$ws.Cells.Item(17, 12).Value = "stimuli/img_2js6m.png"
$ws.Cells.Item(17, 13).Value = 40.02777777777778
$ws.Cells.Item(17, 14).Value = 20.88888888888889
$ws.Cells.Item(17, 15).Value = 30.45833333333334
$ws.Cells.Item(17, 16).Value = 36
$ws.Cells.Item(17, 17).Value = 2
$ws.Cells.Item(17, 18).Value = 2
$ws.Cells.Item(17, 19).Value = 2
$ws.Cells.Item(17, 20).Value = 2
$ws.Cells.Item(17, 21).Value = 2
$ws.Cells.Item(17, 22).Value = 2

# Row 18
$ws.Cells.Item(18, 12).Value = "stimuli/img_a9acb.png"
$ws.Cells.Item(18, 13).Value = 77.11428571428571
$ws.Cells.Item(18, 14).Value = 58.42857142857143
$ws.Cells.Item(18, 15).Value = 67.77142857142857
$ws.Cells.Item(18, 16).Value = 35
$ws.Cells.Item(18, 17).Value = 7
$ws.Cells.Item(18, 18).Value = 7
$ws.Cells.Item(18, 19).Value = 7
$ws.Cells.Item(18, 20).Value = 7
$ws.Cells.Item(18, 21).Value = 7
$ws.Cells.Item(18, 22).Value = 7

# Row 19
$ws.Cells.Item(19, 9).ClearContents()
$ws.Cells.Item(19, 10).Value = "new"
$ws.Cells.Item(19, 11).Value = "f"
$ws.Cells.Item(19, 12).Value = "stimuli/img_zt893.png"
$ws.Cells.Item(19, 13).Value = 68.53191489361703
$ws.Cells.Item(19, 14).Value = 49.19148936170212
$ws.Cells.Item(19, 15).Value = 58.86170212765958
$ws.Cells.Item(19, 16).Value = 47
$ws.Cells.Item(19, 17).Value = 5
$ws.Cells.Item(19, 18).Value = 5
$ws.Cells.Item(19, 19).Value = 5
$ws.Cells.Item(19, 20).Value = 5
$ws.Cells.Item(19, 21).Value = 5
$ws.Cells.Item(19, 22).Value = 5

# Row 20
$ws.Cells.Item(20, 9).ClearContents()
$ws.Cells.Item(20, 10).Value = "new"
$ws.Cells.Item(20, 11).Value = "f"
$ws.Cells.Item(20, 12).Value = "stimuli/img_0eflx.png"
$ws.Cells.Item(20, 13).Value = 76.05128205128206
$ws.Cells.Item(20, 14).Value = 53.53846153846154
$ws.Cells.Item(20, 15).Value = 64.7948717948718
$ws.Cells.Item(20, 17).Value = 6
$ws.Cells.Item(20, 18).Value = 6
$ws.Cells.Item(20, 19).Value = 6
$ws.Cells.Item(20, 20).Value = 6
$ws.Cells.Item(20, 21).Value = 6
$ws.Cells.Item(20, 22).Value = 6

# Row 21
$ws.Cells.Item(21, 12).Value = "stimuli/img_wyctg.png"
$ws.Cells.Item(21, 13).Value = 33.44736842105263
$ws.Cells.Item(21, 14).Value = 11.39473684210526
$ws.Cells.Item(21, 15).Value = 22.42105263157895
$ws.Cells.Item(21, 16).Value = 38
$ws.Cells.Item(21, 17).Value = 1
$ws.Cells.Item(21, 18).Value = 1
$ws.Cells.Item(21, 19).Value = 1
$ws.Cells.Item(21, 20).Value = 1
$ws.Cells.Item(21, 21).Value = 1
$ws.Cells.Item(21, 22).Value = 1

# Row 22
$ws.Cells.Item(22, 8).ClearContents()
$ws.Cells.Item(22, 9).ClearContents()
$ws.Cells.Item(22, 10).Value = "catch"
$ws.Cells.Item(22, 11).Value = "f"
$ws.Cells.Item(22, 12).Value = "stimuli/catch_21.jpg"
$ws.Cells.Item(22, 13).ClearContents()
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(22, 15).ClearContents()
$ws.Cells.Item(22, 16).ClearContents()
$ws.Cells.Item(22, 17).ClearContents()
$ws.Cells.Item(22, 18).ClearContents()
$ws.Cells.Item(22, 19).ClearContents()
$ws.Cells.Item(22, 20).ClearContents()
$ws.Cells.Item(22, 21).ClearContents()
$ws.Cells.Item(22, 22).ClearContents()

# Row 23
$ws.Cells.Item(23, 12).Value = "stimuli/img_juob3.png"
$ws.Cells.Item(23, 13).Value = 79.92105263157895
$ws.Cells.Item(23, 14).Value = 59.78947368421053
$ws.Cells.Item(23, 15).Value = 69.85526315789474
$ws.Cells.Item(23, 16).Value = 38
$ws.Cells.Item(23, 17).Value = 7
$ws.Cells.Item(23, 18).Value = 7
$ws.Cells.Item(23, 19).Value = 7
$ws.Cells.Item(23, 20).Value = 7
$ws.Cells.Item(23, 21).Value = 7
$ws.Cells.Item(23, 22).Value = 7

# Row 24
$ws.Cells.Item(24, 9).Value = "target"
$ws.Cells.Item(24, 10).Value = "old"
$ws.Cells.Item(24, 11).Value = "j"
$ws.Cells.Item(24, 12).Value = "stimuli/img_x0u5z.png"
$ws.Cells.Item(24, 13).Value = 92
$ws.Cells.Item(24, 14).Value = 78.16216216216216
$ws.Cells.Item(24, 15).Value = 85.08108108108108
$ws.Cells.Item(24, 16).Value = 37
$ws.Cells.Item(24, 17).Value = 10
$ws.Cells.Item(24, 18).Value = 10
$ws.Cells.Item(24, 19).Value = 10
$ws.Cells.Item(24, 20).Value = 10
$ws.Cells.Item(24, 21).Value = 10
$ws.Cells.Item(24, 22).Value = 10

# Row 25
$ws.Cells.Item(25, 12).Value = "stimuli/img_gqy6z.png"
$ws.Cells.Item(25, 13).Value = 86.47368421052632
$ws.Cells.Item(25, 14).Value = 68.42105263157895
$ws.Cells.Item(25, 15).Value = 77.44736842105263
$ws.Cells.Item(25, 16).Value = 38
$ws.Cells.Item(25, 17).Value = 9
$ws.Cells.Item(25, 18).Value = 9
$ws.Cells.Item(25, 19).Value = 9
$ws.Cells.Item(25, 20).Value = 9
$ws.Cells.Item(25, 21).Value = 9
$ws.Cells.Item(25, 22).Value = 9

# Row 26
$ws.Cells.Item(26, 12).Value = "stimuli/img_ozxpp.png"
$ws.Cells.Item(26, 13).Value = 26.26470588235294
$ws.Cells.Item(26, 14).Value = 11.47058823529412
$ws.Cells.Item(26, 15).Value = 18.86764705882353
$ws.Cells.Item(26, 16).Value = 34
$ws.Cells.Item(26, 17).Value = 1
$ws.Cells.Item(26, 18).Value = 1
$ws.Cells.Item(26, 19).Value = 1
$ws.Cells.Item(26, 20).Value = 1
$ws.Cells.Item(26, 21).Value = 1
$ws.Cells.Item(26, 22).Value = 1

# Row 27
$ws.Cells.Item(27, 12).Value = "stimuli/img_kljj4.png"
$ws.Cells.Item(27, 13).Value = 64.34999999999999
$ws.Cells.Item(27, 14).Value = 44.15
$ws.Cells.Item(27, 15).Value = 54.25
$ws.Cells.Item(27, 16).Value = 40
$ws.Cells.Item(27, 17).Value = 4
$ws.Cells.Item(27, 18).Value = 4
$ws.Cells.Item(27, 19).Value = 4
$ws.Cells.Item(27, 20).Value = 4
$ws.Cells.Item(27, 21).Value = 4
$ws.Cells.Item(27, 22).Value = 4

# Row 28
$ws.Cells.Item(28, 12).Value = "stimuli/img_bpyv5.png"
$ws.Cells.Item(28, 13).Value = 59.05882352941177
$ws.Cells.Item(28, 14).Value = 37.55882352941177
$ws.Cells.Item(28, 15).Value = 48.30882352941177
$ws.Cells.Item(28, 17).Value = 3
$ws.Cells.Item(28, 18).Value = 3
$ws.Cells.Item(28, 19).Value = 3
$ws.Cells.Item(28, 20).Value = 3
$ws.Cells.Item(28, 21).Value = 3
$ws.Cells.Item(28, 22).Value = 3

# Row 29
$ws.Cells.Item(29, 12).Value = "stimuli/img_u2o6z.png"
$ws.Cells.Item(29, 13).Value = 58.6
$ws.Cells.Item(29, 14).Value = 38.2
$ws.Cells.Item(29, 15).Value = 48.40000000000001
$ws.Cells.Item(29, 16).Value = 30
$ws.Cells.Item(29, 17).Value = 3
$ws.Cells.Item(29, 18).Value = 3
$ws.Cells.Item(29, 19).Value = 3
$ws.Cells.Item(29, 20).Value = 3
$ws.Cells.Item(29, 21).Value = 3
$ws.Cells.Item(29, 22).Value = 3

# Row 30
$ws.Cells.Item(30, 9).Value = "target"
$ws.Cells.Item(30, 10).Value = "old"
$ws.Cells.Item(30, 11).Value = "j"
$ws.Cells.Item(30, 12).Value = "stimuli/img_t2ioc.png"
$ws.Cells.Item(30, 13).Value = 88.18918918918919
$ws.Cells.Item(30, 14).Value = 74.05405405405405
$ws.Cells.Item(30, 15).Value = 81.12162162162161
$ws.Cells.Item(30, 16).Value = 37
$ws.Cells.Item(30, 17).Value = 10
$ws.Cells.Item(30, 18).Value = 10
$ws.Cells.Item(30, 19).Value = 10
$ws.Cells.Item(30, 20).Value = 10
$ws.Cells.Item(30, 21).Value = 10
$ws.Cells.Item(30, 22).Value = 10

# Row 31
$ws.Cells.Item(31, 8).Value = "bedrooms"
$ws.Cells.Item(31, 9).Value = "target"
$ws.Cells.Item(31, 10).Value = "old"
$ws.Cells.Item(31, 11).Value = "j"
$ws.Cells.Item(31, 12).Value = "stimuli/img_e26ut.png"
$ws.Cells.Item(31, 13).Value = 81.07692307692308
$ws.Cells.Item(31, 14).Value = 61.28205128205128
$ws.Cells.Item(31, 15).Value = 71.17948717948718
$ws.Cells.Item(31, 16).Value = 39
$ws.Cells.Item(31, 17).Value = 8
$ws.Cells.Item(31, 18).Value = 8
$ws.Cells.Item(31, 19).Value = 8
$ws.Cells.Item(31, 20).Value = 8
$ws.Cells.Item(31, 21).Value = 8
$ws.Cells.Item(31, 22).Value = 8

# Row 32
$ws.Cells.Item(32, 9).ClearContents()
$ws.Cells.Item(32, 10).Value = "new"
$ws.Cells.Item(32, 11).Value = "f"
$ws.Cells.Item(32, 12).Value = "stimuli/img_ca8kd.png"
$ws.Cells.Item(32, 13).Value = 92.05405405405405
$ws.Cells.Item(32, 14).Value = 73.02702702702703
$ws.Cells.Item(32, 15).Value = 82.54054054054055
$ws.Cells.Item(32, 16).Value = 37
$ws.Cells.Item(32, 17).Value = 10
$ws.Cells.Item(32, 18).Value = 10
$ws.Cells.Item(32, 19).Value = 10
$ws.Cells.Item(32, 20).Value = 10
$ws.Cells.Item(32, 21).Value = 10
$ws.Cells.Item(32, 22).Value = 10

# Row 33
$ws.Cells.Item(33, 12).Value = "stimuli/img_jp28n.png"
$ws.Cells.Item(33, 13).Value = 65.02564102564102
$ws.Cells.Item(33, 14).Value = 44.97435897435897
$ws.Cells.Item(33, 15).Value = 55
$ws.Cells.Item(33, 16).Value = 39
$ws.Cells.Item(33, 17).Value = 4
$ws.Cells.Item(33, 18).Value = 4
$ws.Cells.Item(33, 19).Value = 4
$ws.Cells.Item(33, 20).Value = 4
$ws.Cells.Item(33, 21).Value = 4
$ws.Cells.Item(33, 22).Value = 5

# Row 34
$ws.Cells.Item(34, 9).ClearContents()
$ws.Cells.Item(34, 10).Value = "new"
$ws.Cells.Item(34, 11).Value = "f"
$ws.Cells.Item(34, 12).Value = "stimuli/img_g2akb.png"
$ws.Cells.Item(34, 13).Value = 87.875
$ws.Cells.Item(34, 14).Value = 79
$ws.Cells.Item(34, 15).Value = 83.4375
$ws.Cells.Item(34, 16).Value = 40
$ws.Cells.Item(34, 17).Value = 10
$ws.Cells.Item(34, 18).Value = 10
$ws.Cells.Item(34, 19).Value = 10
$ws.Cells.Item(34, 20).Value = 10
$ws.Cells.Item(34, 21).Value = 10
$ws.Cells.Item(34, 22).Value = 10

# Row 35
$ws.Cells.Item(35, 12).Value = "stimuli/img_oou46.png"
$ws.Cells.Item(35, 13).Value = 75.70270270270271
$ws.Cells.Item(35, 14).Value = 54.86486486486486
$ws.Cells.Item(35, 15).Value = 65.28378378378379
$ws.Cells.Item(35, 17).Value = 6
$ws.Cells.Item(35, 18).Value = 6
$ws.Cells.Item(35, 19).Value = 6
$ws.Cells.Item(35, 20).Value = 6
$ws.Cells.Item(35, 21).Value = 6
$ws.Cells.Item(35, 22).Value = 6

# Row 37
$ws.Cells.Item(37, 12).Value = "stimuli/img_uxxo0.png"
$ws.Cells.Item(37, 13).Value = 71.74418604651163
$ws.Cells.Item(37, 14).Value = 48.44186046511628
$ws.Cells.Item(37, 15).Value = 60.09302325581395
$ws.Cells.Item(37, 16).Value = 43
$ws.Cells.Item(37, 17).Value = 5
$ws.Cells.Item(37, 18).Value = 5
$ws.Cells.Item(37, 19).Value = 5
$ws.Cells.Item(37, 20).Value = 5
$ws.Cells.Item(37, 21).Value = 5
$ws.Cells.Item(37, 22).Value = 5

# Row 38
$ws.Cells.Item(38, 12).Value = "stimuli/img_1vq1v.png"
$ws.Cells.Item(38, 13).Value = 69.42857142857143
$ws.Cells.Item(38, 14).Value = 46.59523809523809
$ws.Cells.Item(38, 15).Value = 58.01190476190476
$ws.Cells.Item(38, 16).Value = 42
$ws.Cells.Item(38, 17).Value = 5
$ws.Cells.Item(38, 18).Value = 5
$ws.Cells.Item(38, 19).Value = 5
$ws.Cells.Item(38, 20).Value = 5
$ws.Cells.Item(38, 21).Value = 5

# Row 39
$ws.Cells.Item(39, 12).Value = "stimuli/img_cogrz.png"
$ws.Cells.Item(39, 13).Value = 60.5
$ws.Cells.Item(39, 14).Value = 39.71428571428572
$ws.Cells.Item(39, 15).Value = 50.10714285714286
$ws.Cells.Item(39, 16).Value = 42
$ws.Cells.Item(39, 17).Value = 3
$ws.Cells.Item(39, 18).Value = 3
$ws.Cells.Item(39, 19).Value = 3
$ws.Cells.Item(39, 20).Value = 3
$ws.Cells.Item(39, 21).Value = 3
$ws.Cells.Item(39, 22).Value = 3

# Row 40
$ws.Cells.Item(40, 12).Value = "stimuli/img_3h4c9.png"
$ws.Cells.Item(40, 13).Value = 85.47619047619048
$ws.Cells.Item(40, 14).Value = 67.26190476190476
$ws.Cells.Item(40, 15).Value = 76.36904761904762
$ws.Cells.Item(40, 16).Value = 42
$ws.Cells.Item(40, 17).Value = 9
$ws.Cells.Item(40, 18).Value = 9
$ws.Cells.Item(40, 19).Value = 9
$ws.Cells.Item(40, 20).Value = 9
$ws.Cells.Item(40, 21).Value = 9
$ws.Cells.Item(40, 22).Value = 9

# Row 41
$ws.Cells.Item(41, 12).Value = "stimuli/img_fnu4h.png"
$ws.Cells.Item(41, 13).Value = 85.87179487179488
$ws.Cells.Item(41, 14).Value = 70.71794871794872
$ws.Cells.Item(41, 15).Value = 78.2948717948718
$ws.Cells.Item(41, 16).Value = 39
$ws.Cells.Item(41, 17).Value = 9
$ws.Cells.Item(41, 18).Value = 9
$ws.Cells.Item(41, 19).Value = 9
$ws.Cells.Item(41, 20).Value = 9
$ws.Cells.Item(41, 21).Value = 9
$ws.Cells.Item(41, 22).Value = 9

# Row 42
$ws.Cells.Item(42, 12).Value = "stimuli/img_c4uwt.png"
$ws.Cells.Item(42, 13).Value = 44.48387096774194
$ws.Cells.Item(42, 14).Value = 30.06451612903226
$ws.Cells.Item(42, 15).Value = 37.2741935483871
$ws.Cells.Item(42, 16).Value = 31
